# Add I0 and IF columns to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from an existing header cell (H1) so the new
# header cells match the existing bold / bordered / centered look.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the two new columns, rows 2-15.
$values = @(
    @(9, 9),
    @(8, 9),
    @(7, 8),
    @(4, 5),
    @(8, 9),
    @(4, 5),
    @(7, 7),
    @(6, 7),
    @(4, 6),
    @(6, 7),
    @(7, 9),
    @(5, 6),
    @(8, 8),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
